$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "26.867.51"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "1.842.11"
$ws.Range("E3").Value = "  +1.57%  "
Set-TextValue "D4" "1.006"
$ws.Range("E4").Value = "  +0.17%  "
Set-TextValue "D5" "309.41"
$ws.Range("E5").Value = "  +1.17%  "
Set-TextValue "D6" "1.005"
$ws.Range("E6").Value = "  +0.20%  "
Set-TextValue "D7" "0.4702"
$ws.Range("E7").Value = "  +3.73%  "
$ws.Range("E8").Value = "  +1.97%  "
Set-TextValue "D9" "0.07144"
$ws.Range("E9").Value = "  +0.58%  "
Set-TextValue "D10" "0.9270"
$ws.Range("E10").Value = "  +4.21%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.930.91"
$ws.Range("E11").Value = "  +6.32%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D12" "19.54"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D13" "0.07693"
$ws.Range("E13").Value = "  -0.72%  "
Set-TextValue "D14" "5.286"
$ws.Range("E14").Value = "  +0.50%  "
Set-TextValue "D15" "6.394"
$ws.Range("E15").Value = "  +1.59%  "
Set-TextValue "D16" "88.15"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("E17").Value = "  +0.38%  "
Set-TextValue "D18" "0.000008623"
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "26.917.48"
$ws.Range("E20").Value = "  +1.29%  "
Set-TextValue "D21" "14.44"
$ws.Range("E21").Value = "  +2.33%  "
Set-TextValue "D22" "5.013"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("E23").Value = "  +1.00%  "
Set-TextValue "D24" "1.933"
$ws.Range("E24").Value = "  -0.42%  "
Set-TextValue "D25" "151.88"
$ws.Range("E25").Value = "  +0.01%  "
Set-TextValue "D26" "18.24"
$ws.Range("E26").Value = "  +2.54%  "
Set-TextValue "D27" "2.008"
$ws.Range("E27").Value = "  -0.53%  "
Set-TextValue "D28" "114.09"
$ws.Range("E28").Value = "  +1.65%  "
Set-TextValue "D29" "4.879"
$ws.Range("E29").Value = "  +1.07%  "
Set-TextValue "D30" "0.08823"
$ws.Range("E30").Value = "  +1.34%  "
Set-TextValue "D31" "3.215"
$ws.Range("E31").Value = "  +2.64%  "
Set-TextValue "D32" "1.179"
$ws.Range("E32").Value = "  +6.32%  "
Set-TextValue "D33" "0.7465"
$ws.Range("E33").Value = "  +0.90%  "
Set-TextValue "D34" "2.781"
$ws.Range("E34").Value = "  +2.69%  "
Set-TextValue "D35" "4.468"
$ws.Range("E35").Value = "  +1.08%  "
Set-TextValue "D36" "1.086"
$ws.Range("E36").Value = "  +1.37%  "
Set-TextValue "D37" "0.01937"
Set-TextValue "D38" "0.05202"
$ws.Range("E38").Value = "  +2.38%  "
$ws.Range("E39").Value = "  +1.46%  "
Set-TextValue "D40" "0.5206"
$ws.Range("E40").Value = "  +2.19%  "
Set-TextValue "D41" "6.960"
$ws.Range("E41").Value = "  +2.79%  "
Set-TextValue "D42" "0.1509"
$ws.Range("E42").Value = "  +0.10%  "
Set-TextValue "D43" "8.148"
$ws.Range("E43").Value = "  +1.49%  "
Set-TextValue "D44" "10.41"
$ws.Range("E44").Value = "  +4.33%  "
Set-TextValue "D45" "0.4694"
$ws.Range("E45").Value = "  +0.00%  "
Set-TextValue "D46" "1.006"
$ws.Range("E46").Value = "  +0.25%  "
Set-TextValue "D47" "101.34"
$ws.Range("E47").Value = "  +2.74%  "
Set-TextValue "D48" "1.595"
$ws.Range("E48").Value = "  +1.90%  "
Set-TextValue "D49" "65.60"
$ws.Range("E49").Value = "  +2.85%  "
Set-TextValue "D50" "0.06035"
$ws.Range("E50").Value = "  +0.95%  "
Set-TextValue "D51" "0.8910"
$ws.Range("E51").Value = "  +5.23%  "
